$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 14
$ws.Range("H3").Value = 14

$ws.Range("F8").Value = 9
$ws.Range("H8").Value = 9

$ws.Range("E9").Value = 23

$ws.Range("E15").Value = 154

$ws.Range("E17").Value = 109
$ws.Range("F17").Value = 49
$ws.Range("H17").Value = 49

$ws.Range("E19").Value = 53

$ws.Range("F32").Value = 5
$ws.Range("H32").Value = 5

$ws.Range("E33").Value = 39

$ws.Range("E34").Value = 19

$ws.Range("E37").Value = 48
$ws.Range("F37").Value = 27
$ws.Range("H37").Value = 27

$ws.Range("E38").Value = 68

$ws.Range("E41").Value = 38

$ws.Range("F42").Value = 16
$ws.Range("H42").Value = 16

$ws.Range("E47").Value = 57

$ws.Range("E66").Value = 33
$ws.Range("F66").Value = 21
$ws.Range("H66").Value = 21

$ws.Range("E68").Value = 15
$ws.Range("F68").Value = 9
$ws.Range("H68").Value = 9

$ws.Range("E70").Value = 42
$ws.Range("F70").Value = 19
$ws.Range("H70").Value = 19

$ws.Range("F72").Value = 18
$ws.Range("H72").Value = 18

$ws.Range("E75").Value = 12

$ws.Range("F78").Value = 19
$ws.Range("H78").Value = 19

$ws.Range("E79").Value = 35
$ws.Range("F79").Value = 15
$ws.Range("H79").Value = 15

$ws.Range("E80").Value = 24
